# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (G1) onto
# the new header cell so the look (bold font, border, centered) matches.
$ws.Range("G1").Copy($ws.Range("H1"))

# Header text.
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
